$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1205.1578
$ws.Range("I70").Value = 990.7273
$ws.Range("K70").Value = 2972.1819
$ws.Range("M70").Value = -2702.1819
$ws.Range("H73").Value = 1205.1578
$ws.Range("I73").Value = 990.7273
$ws.Range("K73").Value = 2972.1819
$ws.Range("M73").Value = -2036.1819
$ws.Range("H108").Value = 31580
$ws.Range("J108").Value = 31580
$ws.Range("L108").Value = 31580
$ws.Range("N108").Value = -39260
$ws.Range("H120").Value = 49736
$ws.Range("J120").Value = 49736
$ws.Range("L120").Value = 49736
$ws.Range("N120").Value = -59412
$ws.Range("H124").Value = 43304.6
$ws.Range("J124").Value = 43304.6
$ws.Range("L124").Value = 43304.6
$ws.Range("N124").Value = -53124.6
$ws.Range("H128").Value = 36971.715
$ws.Range("J128").Value = 36971.715
$ws.Range("L128").Value = 36971.715
$ws.Range("N128").Value = -46931.715
$ws.Range("H130").Value = 43634
$ws.Range("J130").Value = 43634
$ws.Range("L130").Value = 43634
$ws.Range("N130").Value = -53674

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 2000
$ws.Range("I36").Value = 2000
$ws.Range("K36").Value = 2000
$ws.Range("M36").Value = -1654
$ws.Range("H45").Value = 2053.9443
$ws.Range("I45").Value = 1926.5714
$ws.Range("K45").Value = 1926.5714
$ws.Range("M45").Value = -1549.5714
$ws.Range("H61").Value = 2397.0356
$ws.Range("I61").Value = 1393.3334
$ws.Range("J61").Value = 4203.7
$ws.Range("K61").Value = 1393.3334
$ws.Range("L61").Value = 4203.7
$ws.Range("M61").Value = -1181.3334
$ws.Range("N61").Value = -4627.7
$ws.Range("H102").Value = 17619.154
$ws.Range("I102").Value = 1558.7778
$ws.Range("J102").Value = 53755
$ws.Range("K102").Value = 1558.7778
$ws.Range("L102").Value = 53755
$ws.Range("M102").Value = 63.22219999999993
$ws.Range("N102").Value = -56999
$ws.Range("H109").Value = 38874.332
$ws.Range("J109").Value = 38874.332
$ws.Range("L109").Value = 38874.332
$ws.Range("N109").Value = -41648.332
$ws.Range("H113").Value = 40953.6
$ws.Range("J113").Value = 40953.6
$ws.Range("L113").Value = 40953.6
$ws.Range("N113").Value = -49631.6
$ws.Range("H117").Value = 47309.75
$ws.Range("J117").Value = 47309.75
$ws.Range("L117").Value = 47309.75
$ws.Range("N117").Value = -56487.75
$ws.Range("H119").Value = 50694
$ws.Range("J119").Value = 50694
$ws.Range("L119").Value = 50694
$ws.Range("N119").Value = -60370
$ws.Range("H130").Value = 37273.668
$ws.Range("J130").Value = 37273.668
$ws.Range("L130").Value = 37273.668
$ws.Range("N130").Value = -47313.668
$ws.Range("H136").Value = 2397.0356
$ws.Range("I136").Value = 1393.3334
$ws.Range("J136").Value = 4203.7
$ws.Range("K136").Value = 4180.0002
$ws.Range("L136").Value = 12611.1
$ws.Range("M136").Value = -1630.0002
$ws.Range("N136").Value = -17711.1

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 44494
$ws.Range("J112").Value = 44494
$ws.Range("L112").Value = 44494
$ws.Range("N112").Value = -47448
$ws.Range("H124").Value = 49496
$ws.Range("J124").Value = 49496
$ws.Range("L124").Value = 49496
$ws.Range("N124").Value = -59316
$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620
$ws.Range("H126").Value = 50772
$ws.Range("J126").Value = 50772
$ws.Range("L126").Value = 50772
$ws.Range("N126").Value = -60652

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45908.4
$ws.Range("J20").Value = 45908.4
$ws.Range("L20").Value = 45908.4
$ws.Range("N20").Value = -46380.4
$ws.Range("H22").Value = 518.1
$ws.Range("I22").Value = 576.2
$ws.Range("J22").Value = 460
$ws.Range("K22").Value = 576.2
$ws.Range("L22").Value = 460
$ws.Range("M22").Value = -226.2
$ws.Range("N22").Value = -1160
$ws.Range("H30").Value = 45908.4
$ws.Range("J30").Value = 45908.4
$ws.Range("L30").Value = 45908.4
$ws.Range("N30").Value = -46090.4
$ws.Range("H111").Value = 47000
$ws.Range("J111").Value = 47000
$ws.Range("L111").Value = 47000
$ws.Range("N111").Value = -55180
$ws.Range("H112").Value = 31502.2
$ws.Range("J112").Value = 31502.2
$ws.Range("L112").Value = 31502.2
$ws.Range("N112").Value = -34456.2
$ws.Range("H116").Value = 44489.668
$ws.Range("J116").Value = 44489.668
$ws.Range("L116").Value = 44489.668
$ws.Range("N116").Value = -53667.668
$ws.Range("H119").Value = 45486.75
$ws.Range("J119").Value = 45486.75
$ws.Range("L119").Value = 45486.75
$ws.Range("N119").Value = -55162.75
$ws.Range("H128").Value = 45908.4
$ws.Range("J128").Value = 45908.4
$ws.Range("L128").Value = 45908.4
$ws.Range("N128").Value = -55868.4
$ws.Range("H134").Value = 425255.3
$ws.Range("I134").Value = 786.5833
$ws.Range("J134").Value = 1557171.9
$ws.Range("K134").Value = 2359.7499
$ws.Range("L134").Value = 4671515.699999999
$ws.Range("M134").Value = 175.2501000000002
$ws.Range("N134").Value = -4676585.699999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1314.4783
$ws.Range("I102").Value = 1330.8667
$ws.Range("J102").Value = 1283.75
$ws.Range("K102").Value = 1330.8667
$ws.Range("L102").Value = 1283.75
$ws.Range("M102").Value = 291.1333
$ws.Range("N102").Value = -4527.75
$ws.Range("H110").Value = 47749
$ws.Range("J110").Value = 47749
$ws.Range("L110").Value = 47749
$ws.Range("N110").Value = -55929
$ws.Range("H130").Value = 48785.715
$ws.Range("J130").Value = 48785.715
$ws.Range("L130").Value = 48785.715
$ws.Range("N130").Value = -58825.715

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 47325.332
$ws.Range("J108").Value = 47325.332
$ws.Range("L108").Value = 47325.332
$ws.Range("N108").Value = -55005.332
$ws.Range("H112").Value = 31909.666
$ws.Range("J112").Value = 35291.6
$ws.Range("L112").Value = 35291.6
$ws.Range("N112").Value = -38245.6
$ws.Range("H118").Value = 38045.75
$ws.Range("J118").Value = 38045.75
$ws.Range("L118").Value = 38045.75
$ws.Range("N118").Value = -41359.75
$ws.Range("H120").Value = 53108.668
$ws.Range("J120").Value = 53108.668
$ws.Range("L120").Value = 53108.668
$ws.Range("N120").Value = -62784.668
$ws.Range("H121").Value = 42280
$ws.Range("J121").Value = 42280
$ws.Range("L121").Value = 42280
$ws.Range("N121").Value = -45774
$ws.Range("H128").Value = 39992
$ws.Range("J128").Value = 39992
$ws.Range("L128").Value = 39992
$ws.Range("N128").Value = -49952
$ws.Range("H130").Value = 42552.668
$ws.Range("J130").Value = 42552.668
$ws.Range("L130").Value = 42552.668
$ws.Range("N130").Value = -52592.668

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1596.4445
$ws.Range("I96").Value = 675.75
$ws.Range("J96").Value = 2333
$ws.Range("K96").Value = 675.75
$ws.Range("L96").Value = 2333
$ws.Range("M96").Value = 697.25
$ws.Range("N96").Value = -5079
$ws.Range("H110").Value = 46639.332
$ws.Range("J110").Value = 46639.332
$ws.Range("L110").Value = 46639.332
$ws.Range("N110").Value = -54819.332
$ws.Range("H117").Value = 40903.4
$ws.Range("J117").Value = 40903.4
$ws.Range("L117").Value = 40903.4
$ws.Range("N117").Value = -50081.4
$ws.Range("H120").Value = 42100.8
$ws.Range("J120").Value = 42100.8
$ws.Range("L120").Value = 42100.8
$ws.Range("N120").Value = -51776.8
$ws.Range("H121").Value = 43256.668
$ws.Range("J121").Value = 43256.668
$ws.Range("L121").Value = 43256.668
$ws.Range("N121").Value = -46750.668
